$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Paragraph 1: "SELECT * FROM gamebar.employees;" -----------------
# Drop the gramStart/gramEnd proofErr marks, keep everything else as-is.
$p1 = $d.Paragraphs.Item(1)
$p1xml = '<w:p ' + $wns + ' w14:paraId="195FD732" w14:textId="117CF09B" w:rsidR="00AF2A04" w:rsidRDefault="00646D15" w:rsidP="00646D15"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r w:rsidRPr="00646D15"><w:t xml:space="preserve">SELECT * FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00646D15"><w:t>gamebar.employees</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00646D15"><w:t>;</w:t></w:r></w:p>'
$p1.Range.InsertXML($p1xml)

# --- Paragraph 3: "DELETE * FROM gamebar.employees WHERE(ID  = 2);" --
# Same gramStart/gramEnd cleanup, then append the new "manufacturers"
# exercise paragraphs right after it (replaces p3's range with p3 + the
# seven new paragraphs in one shot).
$p3 = $d.Paragraphs.Item(3)
$p3xml = '<w:p ' + $wns + ' w14:paraId="35F8C1E7" w14:textId="7ACD9745" w:rsidR="00646D15" w:rsidRDefault="00646D15" w:rsidP="00646D15"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t>DELETE</w:t></w:r><w:r w:rsidRPr="00646D15"><w:t xml:space="preserve"> * FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00646D15"><w:t>gamebar.employees</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> WHERE(ID  = 2)</w:t></w:r><w:r w:rsidRPr="00646D15"><w:t>;</w:t></w:r></w:p>'

$blankPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr></w:p>'

$insertPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">INSERT INTO manufacturers(name, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>established_on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) VALUE</w:t></w:r></w:p>'

$bmwPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t>(''BMW'', ''1916-03-01''),</w:t></w:r></w:p>'

$teslaPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t>(''Tesla'', ''2003-01-01''),</w:t></w:r></w:p>'

$ladaPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t>(''Lada'', ''1966-05-01'');</w:t></w:r></w:p>'

$selectPara = '<w:p ' + $wns + '><w:pPr><w:tabs><w:tab w:val="left" w:pos="1517"/></w:tabs></w:pPr><w:r><w:t>SELECT * FROM manufacturers;</w:t></w:r></w:p>'

$fullXml = $p3xml + $blankPara + $insertPara + $bmwPara + $teslaPara + $ladaPara + $blankPara + $selectPara

$p3.Range.InsertXML($fullXml)
